# regen save_data to use K (strikeouts) instead of Strike# proxy for Kyle Finnegan's
# 2024 game log: overwrite the "K" column (G2:G72) with the re-derived strikeout
# counts calculated from the boxscore re-scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values, one per game row (r=2 .. r=72), in order.
$kVals = @(1,0,2,0,1,0,1,2,2,1,0,2,0,0,2,1,0,1,0,0,1,1,0,0,2,1,0,0,1,1,2,2,0,2,1,0,1,3,1,2,1,0,1,3,1,0,1,1,1,0,3,1,1,1,0,1,1,1,1,0,1,2,0,1,1,1,2,2,1,3,1)

$startRow = 2
for ($i = 0; $i -lt $kVals.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kVals[$i]
}
